$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.048758552317011095,
    -0.0099999992874053589,
    -0.0089999992983837984,
    0.28398630982225725,
    -0.0059999993233184057,
    -0.0059999993052386458,
    0.0072314355397384844,
    -0.019999999181043115,
    -0.0059999992982282535,
    -0.0059999992966695004,
    -0.0044999993093668422,
    0.072096089089378701,
    -0.0059999992787975742,
    -0.011999999221778523,
    -0.0059999992682442382,
    -0.0059999992629871102,
    -0.0059999992541417413,
    -0.032733907794731287,
    -0.0089999992940219542,
    -0.0089999992754243863,
    -0.0089999992715963373,
    -0.053714722284969429,
    -0.0089999992922820127,
    -0.041999998998242916,
    -0.041999998992919174,
    -0.0059999993025243725,
    -0.005999999298892611,
    -0.0059999992834347537,
    -0.011999999221968594,
    -0.019999999148831105,
    -0.014999999184494683,
    -0.020999999131974256,
    -0.0059999992600223706
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
